$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Ligand/Receptor-expressing cell counts and all dependent metrics
# (following Dr Hou advice) for rows 2-37, columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T.

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 7).Value = 4.495339
$ws.Cells.Item(2, 8).Value = 8.990677999999999
$ws.Cells.Item(2, 9).Value = 0.1274106381683396
$ws.Cells.Item(2, 10).Value = 0.09709078547908702
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 13).Value = 69.292158
$ws.Cells.Item(2, 14).Value = 138.584316
$ws.Cells.Item(2, 15).Value = 0.2302847784735656
$ws.Cells.Item(2, 16).Value = 0.1742377818809482
$ws.Cells.Item(2, 17).Value = 311.491740251562
$ws.Cells.Item(2, 18).Value = 1245.966961006248
$ws.Cells.Item(2, 19).Value = 0.02934073058577169
$ws.Cells.Item(2, 20).Value = 0.0169168831029551

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 7).Value = 4.495339
$ws.Cells.Item(3, 8).Value = 8.990677999999999
$ws.Cells.Item(3, 9).Value = 0.1274106381683396
$ws.Cells.Item(3, 10).Value = 0.09709078547908702
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 35.08748833333333
$ws.Cells.Item(3, 14).Value = 105.262465
$ws.Cells.Item(3, 15).Value = 0.1166093640500485
$ws.Cells.Item(3, 16).Value = 0.1323432473911474
$ws.Cells.Item(3, 17).Value = 157.7301547168783
$ws.Cells.Item(3, 18).Value = 946.3809283012699
$ws.Cells.Item(3, 19).Value = 0.01485727349002092
$ws.Cells.Item(3, 20).Value = 0.01284930984205964

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 7).Value = 4.495339
$ws.Cells.Item(4, 8).Value = 8.990677999999999
$ws.Cells.Item(4, 9).Value = 0.1274106381683396
$ws.Cells.Item(4, 10).Value = 0.09709078547908702
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 64.44625533333334
$ws.Cells.Item(4, 14).Value = 193.338766
$ws.Cells.Item(4, 15).Value = 0.214179960059658
$ws.Cells.Item(4, 16).Value = 0.2430788613874581
$ws.Cells.Item(4, 17).Value = 289.7077650038913
$ws.Cells.Item(4, 18).Value = 1738.246590023348
$ws.Cells.Item(4, 19).Value = 0.0272888053940705
$ws.Cells.Item(4, 20).Value = 0.02360071758547042

# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 7).Value = 4.495339
$ws.Cells.Item(5, 8).Value = 8.990677999999999
$ws.Cells.Item(5, 9).Value = 0.1274106381683396
$ws.Cells.Item(5, 10).Value = 0.09709078547908702
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 56.56647866666666
$ws.Cells.Item(5, 14).Value = 169.699436
$ws.Cells.Item(5, 15).Value = 0.1879923989202791
$ws.Cells.Item(5, 16).Value = 0.2133578616146428
$ws.Cells.Item(5, 17).Value = 254.2854976429346
$ws.Cells.Item(5, 18).Value = 1525.712985857608
$ws.Cells.Item(5, 19).Value = 0.02395223151722982
$ws.Cells.Item(5, 20).Value = 0.02071508237230402

# Row 6
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 7).Value = 4.495339
$ws.Cells.Item(6, 8).Value = 8.990677999999999
$ws.Cells.Item(6, 9).Value = 0.1274106381683396
$ws.Cells.Item(6, 10).Value = 0.09709078547908702
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 37.47905366666667
$ws.Cells.Item(6, 14).Value = 112.437161
$ws.Cells.Item(6, 15).Value = 0.1245574653776436
$ws.Cells.Item(6, 16).Value = 0.1413637711617457
$ws.Cells.Item(6, 17).Value = 168.4810516308596
$ws.Cells.Item(6, 18).Value = 1010.886309785158
$ws.Cells.Item(6, 19).Value = 0.01586994615239643
$ws.Cells.Item(6, 20).Value = 0.0137251195803798

# Row 7
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 7).Value = 4.495339
$ws.Cells.Item(7, 8).Value = 8.990677999999999
$ws.Cells.Item(7, 9).Value = 0.1274106381683396
$ws.Cells.Item(7, 10).Value = 0.09709078547908702
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 13).Value = 38.0262565
$ws.Cells.Item(7, 14).Value = 76.052513
$ws.Cells.Item(7, 15).Value = 0.1263760331188052
$ws.Cells.Item(7, 16).Value = 0.0956184765640578
$ws.Cells.Item(7, 17).Value = 170.9409138684535
$ws.Cells.Item(7, 18).Value = 683.763655473814
$ws.Cells.Item(7, 19).Value = 0.01610165102885019
$ws.Cells.Item(7, 20).Value = 0.009283672995918046

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 0.4630963333333333
$ws.Cells.Item(8, 8).Value = 1.389289
$ws.Cells.Item(8, 9).Value = 0.01312546158663854
$ws.Cells.Item(8, 10).Value = 0.01500300202804008
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 13).Value = 69.292158
$ws.Cells.Item(8, 14).Value = 138.584316
$ws.Cells.Item(8, 15).Value = 0.2302847784735656
$ws.Cells.Item(8, 16).Value = 0.1742377818809482
$ws.Cells.Item(8, 17).Value = 32.088944298554
$ws.Cells.Item(8, 18).Value = 192.533665791324
$ws.Cells.Item(8, 19).Value = 0.00302259401384235
$ws.Cells.Item(8, 20).Value = 0.002614089794921072

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 0.4630963333333333
$ws.Cells.Item(9, 8).Value = 1.389289
$ws.Cells.Item(9, 9).Value = 0.01312546158663854
$ws.Cells.Item(9, 10).Value = 0.01500300202804008
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 35.08748833333333
$ws.Cells.Item(9, 14).Value = 105.262465
$ws.Cells.Item(9, 15).Value = 0.1166093640500485
$ws.Cells.Item(9, 16).Value = 0.1323432473911474
$ws.Cells.Item(9, 17).Value = 16.24888719304278
$ws.Cells.Item(9, 18).Value = 146.239984737385
$ws.Cells.Item(9, 19).Value = 0.001530551728481261
$ws.Cells.Item(9, 20).Value = 0.001985546009006795

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.4630963333333333
$ws.Cells.Item(10, 8).Value = 1.389289
$ws.Cells.Item(10, 9).Value = 0.01312546158663854
$ws.Cells.Item(10, 10).Value = 0.01500300202804008
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 64.44625533333334
$ws.Cells.Item(10, 14).Value = 193.338766
$ws.Cells.Item(10, 15).Value = 0.214179960059658
$ws.Cells.Item(10, 16).Value = 0.2430788613874581
$ws.Cells.Item(10, 17).Value = 29.84482454193045
$ws.Cells.Item(10, 18).Value = 268.603420877374
$ws.Cells.Item(10, 19).Value = 0.002811210838390818
$ws.Cells.Item(10, 20).Value = 0.003646912650369708

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.4630963333333333
$ws.Cells.Item(11, 8).Value = 1.389289
$ws.Cells.Item(11, 9).Value = 0.01312546158663854
$ws.Cells.Item(11, 10).Value = 0.01500300202804008
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 56.56647866666666
$ws.Cells.Item(11, 14).Value = 169.699436
$ws.Cells.Item(11, 15).Value = 0.1879923989202791
$ws.Cells.Item(11, 16).Value = 0.2133578616146428
$ws.Cells.Item(11, 17).Value = 26.19572886011155
$ws.Cells.Item(11, 18).Value = 235.761559741004
$ws.Cells.Item(11, 19).Value = 0.002467487010608151
$ws.Cells.Item(11, 20).Value = 0.003201008430502781

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.4630963333333333
$ws.Cells.Item(12, 8).Value = 1.389289
$ws.Cells.Item(12, 9).Value = 0.01312546158663854
$ws.Cells.Item(12, 10).Value = 0.01500300202804008
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 37.47905366666667
$ws.Cells.Item(12, 14).Value = 112.437161
$ws.Cells.Item(12, 15).Value = 0.1245574653776436
$ws.Cells.Item(12, 16).Value = 0.1413637711617457
$ws.Cells.Item(12, 17).Value = 17.35641232983656
$ws.Cells.Item(12, 18).Value = 156.207710968529
$ws.Cells.Item(12, 19).Value = 0.001634874227143321
$ws.Cells.Item(12, 20).Value = 0.002120880945431065

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.4630963333333333
$ws.Cells.Item(13, 8).Value = 1.389289
$ws.Cells.Item(13, 9).Value = 0.01312546158663854
$ws.Cells.Item(13, 10).Value = 0.01500300202804008
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 13).Value = 38.0262565
$ws.Cells.Item(13, 14).Value = 76.052513
$ws.Cells.Item(13, 15).Value = 0.1263760331188052
$ws.Cells.Item(13, 16).Value = 0.0956184765640578
$ws.Cells.Item(13, 17).Value = 17.60981995554284
$ws.Cells.Item(13, 18).Value = 105.658919733257
$ws.Cells.Item(13, 19).Value = 0.001658743768172638
$ws.Cells.Item(13, 20).Value = 0.001434564197808662

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 6.798299666666666
$ws.Cells.Item(14, 8).Value = 20.394899
$ws.Cells.Item(14, 9).Value = 0.1926830655017586
$ws.Cells.Item(14, 10).Value = 0.2202455436260365
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 13).Value = 69.292158
$ws.Cells.Item(14, 14).Value = 138.584316
$ws.Cells.Item(14, 15).Value = 0.2302847784735656
$ws.Cells.Item(14, 16).Value = 0.1742377818809482
$ws.Cells.Item(14, 17).Value = 471.068854634014
$ws.Cells.Item(14, 18).Value = 2826.413127804084
$ws.Cells.Item(14, 19).Value = 0.04437197705468001
$ws.Cells.Item(14, 20).Value = 0.03837509499056421

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 6.798299666666666
$ws.Cells.Item(15, 8).Value = 20.394899
$ws.Cells.Item(15, 9).Value = 0.1926830655017586
$ws.Cells.Item(15, 10).Value = 0.2202455436260365
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 35.08748833333333
$ws.Cells.Item(15, 14).Value = 105.262465
$ws.Cells.Item(15, 15).Value = 0.1166093640500485
$ws.Cells.Item(15, 16).Value = 0.1323432473911474
$ws.Cells.Item(15, 17).Value = 238.5352602406705
$ws.Cells.Item(15, 18).Value = 2146.817342166035
$ws.Cells.Item(15, 19).Value = 0.02246864973137392
$ws.Cells.Item(15, 20).Value = 0.0291480104668983

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 6.798299666666666
$ws.Cells.Item(16, 8).Value = 20.394899
$ws.Cells.Item(16, 9).Value = 0.1926830655017586
$ws.Cells.Item(16, 10).Value = 0.2202455436260365
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 64.44625533333334
$ws.Cells.Item(16, 14).Value = 193.338766
$ws.Cells.Item(16, 15).Value = 0.214179960059658
$ws.Cells.Item(16, 16).Value = 0.2430788613874581
$ws.Cells.Item(16, 17).Value = 438.1249561505149
$ws.Cells.Item(16, 18).Value = 3943.124605354634
$ws.Cells.Item(16, 19).Value = 0.04126885127333913
$ws.Cells.Item(16, 20).Value = 0.05353703597027868

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 6.798299666666666
$ws.Cells.Item(17, 8).Value = 20.394899
$ws.Cells.Item(17, 9).Value = 0.1926830655017586
$ws.Cells.Item(17, 10).Value = 0.2202455436260365
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 56.56647866666666
$ws.Cells.Item(17, 14).Value = 169.699436
$ws.Cells.Item(17, 15).Value = 0.1879923989202791
$ws.Cells.Item(17, 16).Value = 0.2133578616146428
$ws.Cells.Item(17, 17).Value = 384.5558730641071
$ws.Cells.Item(17, 18).Value = 3461.002857576964
$ws.Cells.Item(17, 19).Value = 0.03622295171498886
$ws.Cells.Item(17, 20).Value = 0.04699111821820567

# Row 18
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 6.798299666666666
$ws.Cells.Item(18, 8).Value = 20.394899
$ws.Cells.Item(18, 9).Value = 0.1926830655017586
$ws.Cells.Item(18, 10).Value = 0.2202455436260365
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 37.47905366666667
$ws.Cells.Item(18, 14).Value = 112.437161
$ws.Cells.Item(18, 15).Value = 0.1245574653776436
$ws.Cells.Item(18, 16).Value = 0.1413637711617457
$ws.Cells.Item(18, 17).Value = 254.7938380490821
$ws.Cells.Item(18, 18).Value = 2293.144542441739
$ws.Cells.Item(18, 19).Value = 0.02400011426009353
$ws.Cells.Item(18, 20).Value = 0.03113474062854531

# Row 19
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 6.798299666666666
$ws.Cells.Item(19, 8).Value = 20.394899
$ws.Cells.Item(19, 9).Value = 0.1926830655017586
$ws.Cells.Item(19, 10).Value = 0.2202455436260365
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 13).Value = 38.0262565
$ws.Cells.Item(19, 14).Value = 76.052513
$ws.Cells.Item(19, 15).Value = 0.1263760331188052
$ws.Cells.Item(19, 16).Value = 0.0956184765640578
$ws.Cells.Item(19, 17).Value = 258.5138868885311
$ws.Cells.Item(19, 18).Value = 1551.083321331187
$ws.Cells.Item(19, 19).Value = 0.02435052146728316
$ws.Cells.Item(19, 20).Value = 0.02105954335154434

# Row 20
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 8.495531999999999
$ws.Cells.Item(20, 8).Value = 25.486596
$ws.Cells.Item(20, 9).Value = 0.2407874364312792
$ws.Cells.Item(20, 10).Value = 0.275231036505607
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 13).Value = 69.292158
$ws.Cells.Item(20, 14).Value = 138.584316
$ws.Cells.Item(20, 15).Value = 0.2302847784735656
$ws.Cells.Item(20, 16).Value = 0.1742377818809482
$ws.Cells.Item(20, 17).Value = 588.673745638056
$ws.Cells.Item(20, 18).Value = 3532.042473828336
$ws.Cells.Item(20, 19).Value = 0.05544968145779488
$ws.Cells.Item(20, 20).Value = 0.04795564530553124

# Row 21
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 8.495531999999999
$ws.Cells.Item(21, 8).Value = 25.486596
$ws.Cells.Item(21, 9).Value = 0.2407874364312792
$ws.Cells.Item(21, 10).Value = 0.275231036505607
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 35.08748833333333
$ws.Cells.Item(21, 14).Value = 105.262465
$ws.Cells.Item(21, 15).Value = 0.1166093640500485
$ws.Cells.Item(21, 16).Value = 0.1323432473911474
$ws.Cells.Item(21, 17).Value = 298.0868799354599
$ws.Cells.Item(21, 18).Value = 2682.78191941914
$ws.Cells.Item(21, 19).Value = 0.02807806983349295
$ws.Cells.Item(21, 20).Value = 0.03642496915398347

# Row 22
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 8.495531999999999
$ws.Cells.Item(22, 8).Value = 25.486596
$ws.Cells.Item(22, 9).Value = 0.2407874364312792
$ws.Cells.Item(22, 10).Value = 0.275231036505607
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 64.44625533333334
$ws.Cells.Item(22, 14).Value = 193.338766
$ws.Cells.Item(22, 15).Value = 0.214179960059658
$ws.Cells.Item(22, 16).Value = 0.2430788613874581
$ws.Cells.Item(22, 17).Value = 547.505224464504
$ws.Cells.Item(22, 18).Value = 4927.547020180536
$ws.Cells.Item(22, 19).Value = 0.05157184351771881
$ws.Cells.Item(22, 20).Value = 0.06690284697227286

# Row 23
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 8.495531999999999
$ws.Cells.Item(23, 8).Value = 25.486596
$ws.Cells.Item(23, 9).Value = 0.2407874364312792
$ws.Cells.Item(23, 10).Value = 0.275231036505607
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 56.56647866666666
$ws.Cells.Item(23, 14).Value = 169.699436
$ws.Cells.Item(23, 15).Value = 0.1879923989202791
$ws.Cells.Item(23, 16).Value = 0.2133578616146428
$ws.Cells.Item(23, 17).Value = 480.5623296399839
$ws.Cells.Item(23, 18).Value = 4325.060966759856
$ws.Cells.Item(23, 19).Value = 0.04526620780458037
$ws.Cells.Item(23, 20).Value = 0.058722705398818

# Row 24
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 8.495531999999999
$ws.Cells.Item(24, 8).Value = 25.486596
$ws.Cells.Item(24, 9).Value = 0.2407874364312792
$ws.Cells.Item(24, 10).Value = 0.275231036505607
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 37.47905366666667
$ws.Cells.Item(24, 14).Value = 112.437161
$ws.Cells.Item(24, 15).Value = 0.1245574653776436
$ws.Cells.Item(24, 16).Value = 0.1413637711617457
$ws.Cells.Item(24, 17).Value = 318.4044997548839
$ws.Cells.Item(24, 18).Value = 2865.640497793956
$ws.Cells.Item(24, 19).Value = 0.02999187277666061
$ws.Cells.Item(24, 20).Value = 0.03890769726118871

# Row 25
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 8.495531999999999
$ws.Cells.Item(25, 8).Value = 25.486596
$ws.Cells.Item(25, 9).Value = 0.2407874364312792
$ws.Cells.Item(25, 10).Value = 0.275231036505607
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 13).Value = 38.0262565
$ws.Cells.Item(25, 14).Value = 76.052513
$ws.Cells.Item(25, 15).Value = 0.1263760331188052
$ws.Cells.Item(25, 16).Value = 0.0956184765640578
$ws.Cells.Item(25, 17).Value = 323.053278935958
$ws.Cells.Item(25, 18).Value = 1938.319673615748
$ws.Cells.Item(25, 19).Value = 0.03042976104103154
$ws.Cells.Item(25, 20).Value = 0.02631717241381272

# Row 26
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 6.279227333333334
$ws.Cells.Item(26, 8).Value = 18.837682
$ws.Cells.Item(26, 9).Value = 0.1779710855497397
$ws.Cells.Item(26, 10).Value = 0.2034290786507157
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 69.292158
$ws.Cells.Item(26, 14).Value = 138.584316
$ws.Cells.Item(26, 15).Value = 0.2302847784735656
$ws.Cells.Item(26, 16).Value = 0.1742377818809482
$ws.Cells.Item(26, 17).Value = 435.101212499252
$ws.Cells.Item(26, 18).Value = 2610.607274995512
$ws.Cells.Item(26, 19).Value = 0.04098403201052179
$ws.Cells.Item(26, 20).Value = 0.03544503143418566

# Row 27
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 6.279227333333334
$ws.Cells.Item(27, 8).Value = 18.837682
$ws.Cells.Item(27, 9).Value = 0.1779710855497397
$ws.Cells.Item(27, 10).Value = 0.2034290786507157
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 35.08748833333333
$ws.Cells.Item(27, 14).Value = 105.262465
$ws.Cells.Item(27, 15).Value = 0.1166093640500485
$ws.Cells.Item(27, 16).Value = 0.1323432473911474
$ws.Cells.Item(27, 17).Value = 220.3223158006811
$ws.Cells.Item(27, 18).Value = 1982.90084220613
$ws.Cells.Item(27, 19).Value = 0.02075309510525193
$ws.Cells.Item(27, 20).Value = 0.02692246488242485

# Row 28
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 6.279227333333334
$ws.Cells.Item(28, 8).Value = 18.837682
$ws.Cells.Item(28, 9).Value = 0.1779710855497397
$ws.Cells.Item(28, 10).Value = 0.2034290786507157
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 64.44625533333334
$ws.Cells.Item(28, 14).Value = 193.338766
$ws.Cells.Item(28, 15).Value = 0.214179960059658
$ws.Cells.Item(28, 16).Value = 0.2430788613874581
$ws.Cells.Item(28, 17).Value = 404.6726880200459
$ws.Cells.Item(28, 18).Value = 3642.054192180412
$ws.Cells.Item(28, 19).Value = 0.03811783999481722
$ws.Cells.Item(28, 20).Value = 0.04944930881151563

# Row 29
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 6.279227333333334
$ws.Cells.Item(29, 8).Value = 18.837682
$ws.Cells.Item(29, 9).Value = 0.1779710855497397
$ws.Cells.Item(29, 10).Value = 0.2034290786507157
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 56.56647866666666
$ws.Cells.Item(29, 14).Value = 169.699436
$ws.Cells.Item(29, 15).Value = 0.1879923989202791
$ws.Cells.Item(29, 16).Value = 0.2133578616146428
$ws.Cells.Item(29, 17).Value = 355.1937789941503
$ws.Cells.Item(29, 18).Value = 3196.744010947352
$ws.Cells.Item(29, 19).Value = 0.03345721131094177
$ws.Cells.Item(29, 20).Value = 0.04340319321115368

# Row 30
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 6.279227333333334
$ws.Cells.Item(30, 8).Value = 18.837682
$ws.Cells.Item(30, 9).Value = 0.1779710855497397
$ws.Cells.Item(30, 10).Value = 0.2034290786507157
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 37.47905366666667
$ws.Cells.Item(30, 14).Value = 112.437161
$ws.Cells.Item(30, 15).Value = 0.1245574653776436
$ws.Cells.Item(30, 16).Value = 0.1413637711617457
$ws.Cells.Item(30, 17).Value = 235.3394982112002
$ws.Cells.Item(30, 18).Value = 2118.055483900802
$ws.Cells.Item(30, 19).Value = 0.02216762732658334
$ws.Cells.Item(30, 20).Value = 0.02875750172202455

# Row 31
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 6.279227333333334
$ws.Cells.Item(31, 8).Value = 18.837682
$ws.Cells.Item(31, 9).Value = 0.1779710855497397
$ws.Cells.Item(31, 10).Value = 0.2034290786507157
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 13).Value = 38.0262565
$ws.Cells.Item(31, 14).Value = 76.052513
$ws.Cells.Item(31, 15).Value = 0.1263760331188052
$ws.Cells.Item(31, 16).Value = 0.0956184765640578
$ws.Cells.Item(31, 17).Value = 238.7755091991444
$ws.Cells.Item(31, 18).Value = 1432.653055194866
$ws.Cells.Item(31, 19).Value = 0.02249127980162362
$ws.Cells.Item(31, 20).Value = 0.01945157858941133

# Row 32
$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 7).Value = 8.750795
$ws.Cells.Item(32, 8).Value = 17.50159
$ws.Cells.Item(32, 9).Value = 0.2480223127622444
$ws.Cells.Item(32, 10).Value = 0.1890005537105138
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 13).Value = 69.292158
$ws.Cells.Item(32, 14).Value = 138.584316
$ws.Cells.Item(32, 15).Value = 0.2302847784735656
$ws.Cells.Item(32, 16).Value = 0.1742377818809482
$ws.Cells.Item(32, 17).Value = 606.36146976561
$ws.Cells.Item(32, 18).Value = 2425.44587906244
$ws.Cells.Item(32, 19).Value = 0.05711576335095486
$ws.Cells.Item(32, 20).Value = 0.03293103725279094

# Row 33
$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 7).Value = 8.750795
$ws.Cells.Item(33, 8).Value = 17.50159
$ws.Cells.Item(33, 9).Value = 0.2480223127622444
$ws.Cells.Item(33, 10).Value = 0.1890005537105138
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 35.08748833333333
$ws.Cells.Item(33, 14).Value = 105.262465
$ws.Cells.Item(33, 15).Value = 0.1166093640500485
$ws.Cells.Item(33, 16).Value = 0.1323432473911474
$ws.Cells.Item(33, 17).Value = 307.0434174698917
$ws.Cells.Item(33, 18).Value = 1842.26050481935
$ws.Cells.Item(33, 19).Value = 0.02892172416142756
$ws.Cells.Item(33, 20).Value = 0.02501294703677437

# Row 34
$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 7).Value = 8.750795
$ws.Cells.Item(34, 8).Value = 17.50159
$ws.Cells.Item(34, 9).Value = 0.2480223127622444
$ws.Cells.Item(34, 10).Value = 0.1890005537105138
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 64.44625533333334
$ws.Cells.Item(34, 14).Value = 193.338766
$ws.Cells.Item(34, 15).Value = 0.214179960059658
$ws.Cells.Item(34, 16).Value = 0.2430788613874581
$ws.Cells.Item(34, 17).Value = 563.9559689396567
$ws.Cells.Item(34, 18).Value = 3383.73581363794
$ws.Cells.Item(34, 19).Value = 0.05312140904132152
$ws.Cells.Item(34, 20).Value = 0.04594203939755082

# Row 35
$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 7).Value = 8.750795
$ws.Cells.Item(35, 8).Value = 17.50159
$ws.Cells.Item(35, 9).Value = 0.2480223127622444
$ws.Cells.Item(35, 10).Value = 0.1890005537105138
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 56.56647866666666
$ws.Cells.Item(35, 14).Value = 169.699436
$ws.Cells.Item(35, 15).Value = 0.1879923989202791
$ws.Cells.Item(35, 16).Value = 0.2133578616146428
$ws.Cells.Item(35, 17).Value = 495.0016586838733
$ws.Cells.Item(35, 18).Value = 2970.00995210324
$ws.Cells.Item(35, 19).Value = 0.04662630956193007
$ws.Cells.Item(35, 20).Value = 0.04032475398365867

# Row 36
$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 7).Value = 8.750795
$ws.Cells.Item(36, 8).Value = 17.50159
$ws.Cells.Item(36, 9).Value = 0.2480223127622444
$ws.Cells.Item(36, 10).Value = 0.1890005537105138
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 37.47905366666667
$ws.Cells.Item(36, 14).Value = 112.437161
$ws.Cells.Item(36, 15).Value = 0.1245574653776436
$ws.Cells.Item(36, 16).Value = 0.1413637711617457
$ws.Cells.Item(36, 17).Value = 327.9715154309984
$ws.Cells.Item(36, 18).Value = 1967.82909258599
$ws.Cells.Item(36, 19).Value = 0.03089303063476635
$ws.Cells.Item(36, 20).Value = 0.0267178310241763

# Row 37
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 7).Value = 8.750795
$ws.Cells.Item(37, 8).Value = 17.50159
$ws.Cells.Item(37, 9).Value = 0.2480223127622444
$ws.Cells.Item(37, 10).Value = 0.1890005537105138
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 13).Value = 38.0262565
$ws.Cells.Item(37, 14).Value = 76.052513
$ws.Cells.Item(37, 15).Value = 0.1263760331188052
$ws.Cells.Item(37, 16).Value = 0.0956184765640578
$ws.Cells.Item(37, 17).Value = 332.7599752489175
$ws.Cells.Item(37, 18).Value = 1331.03990099567
$ws.Cells.Item(37, 19).Value = 0.03134407601184407
$ws.Cells.Item(37, 20).Value = 0.01807194501556271
